# Correct PRJ unique: false
# The "project_code" uniqueness check was wrong, so the validation-error rows
# it produced ("Duplicate found: PRJ-...") are no longer valid and must be
# removed from the report. Those were rows 7, 9, 11 and 13. Deleting them
# shifts the remaining rows up, shrinking the used range from A1:F13 to A1:F9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete from bottom to top so row indices of not-yet-deleted rows stay valid.
$ws.Rows.Item(13).Delete()
$ws.Rows.Item(11).Delete()
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(7).Delete()
